# SCRUM-111 Updated WebApplication test cases
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (TC-WEB-001): Actual result now mirrors the expected result text,
#     Status moves from "Pending" to "Accepted" (green) ---
$ws.Range("H2").Value = "The homepage of the web application should load successfully without any errors."
$ws.Range("H2").HorizontalAlignment = -4131
$ws.Range("H2").VerticalAlignment = -4108
$ws.Range("H2").WrapText = $true

$ws.Range("I2").Value = "Accepted"
$ws.Range("I2").Interior.Color = 5296274
$ws.Range("I2").HorizontalAlignment = -4108
$ws.Range("I2").VerticalAlignment = -4108
$ws.Range("I2").WrapText = $true

# --- Row 3 (TC-WEB-002): Status moves from "Pending" to "Not Executed" (yellow) ---
$ws.Range("I3").Value = "Not Executed"
$ws.Range("I3").Interior.Color = 65535
$ws.Range("I3").HorizontalAlignment = -4108
$ws.Range("I3").VerticalAlignment = -4108
$ws.Range("I3").WrapText = $true

# --- Row heights re-flowed to fit the new / rewrapped content ---
$ws.Rows(1).RowHeight = 37.5
$ws.Rows(2).RowHeight = 126
$ws.Rows(3).RowHeight = 126

# --- Widen the "Actual result" column so the longer text is readable ---
$ws.Columns("H").ColumnWidth = 23.166666666666668
